# Update for new docker build method
$p = $ppt.ActivePresentation

# --- Slide 13: "Notebook exercise (learn FDTD)" ---
# "Click on work, click on " -> "Click on " (lvl 2 bullet, keeps the
# following run "fdtd_intro.ipynb" untouched).
$s13 = $p.Slides.Item(13)
$tr13 = $s13.Shapes.Item(2).TextFrame.TextRange
$tr13.Paragraphs(3).Characters(1, 24).Text = "Click on "

# --- Slide 14: "Exercise (fdtd1d.c)" ---
# Same fix as slide 13 (keeps the following run "fdtd_intro.c" untouched).
$s14 = $p.Slides.Item(14)
$tr14 = $s14.Shapes.Item(2).TextFrame.TextRange
$tr14.Paragraphs(3).Characters(1, 24).Text = "Click on "

# --- Slide 16: "Exercise (run FDTD in C)" ---
# Remove the "cd work" bullet (the docker image now starts already in
# the work dir) and bump the placeholder's shrink-to-fit font scale.
$s16 = $p.Slides.Item(16)
$tf16 = $s16.Shapes.Item(2).TextFrame
$tr16 = $tf16.TextRange
$tr16.Paragraphs(4).Delete()
$tf16.AutofitFontScale = 0.925

# --- Slide 17: "Exercise (modify fdtd1d.c for glass)" ---
# Drop the standalone "vi discussion" bullet and add two new sub-bullets
# describing the two ways to edit the file.
$s17 = $p.Slides.Item(17)
$tr17 = $s17.Shapes.Item(2).TextFrame.TextRange
$tr17.Paragraphs(1).Delete()
$para17_1 = $tr17.Paragraphs(1)
[void]$para17_1.InsertAfter("`rEdit in binder`rOr use vi from the terminal")
$tr17.Paragraphs(2).IndentLevel = 2
$tr17.Paragraphs(3).IndentLevel = 2
